$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 29

# Column A holds the date as plain text (matching the existing rows' style),
# so force text formatting before assigning the value to avoid Excel
# auto-converting the "MM/DD/YYYY" string into a date serial number.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "09/30/2025"
# Reset the cell style back to Normal/General so this new row doesn't pick
# up a stray numeric-format style that the other data rows don't have.
$cellA.Style = "Normal"

$ws.Cells.Item($row, 2).Value = 0.1412270494756673
$ws.Cells.Item($row, 3).Value = 0.8587729505243327
